$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.479.92'
$ws.Range('D3').Value = '1.813.47'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.64'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.85%  '
$ws.Range('E6').Value = '  +2.88%  '
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '38.36'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.74%  '
$ws.Range('E9').Value = '  -4.11%  '
$ws.Range('E10').Value = '  -2.70%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0973'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.78%  '
$ws.Range('D12').Value = '2.075.04'
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.21'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.33%  '
$ws.Range('D14').Value = '1.815.55'
$ws.Range('E14').Value = '  +0.53%  '
$ws.Range('E15').Value = '  -1.84%  '
$ws.Range('D16').Value = '34.469.21'
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.43'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.10%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.28'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '242.78'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.32%  '
$ws.Range('D20').Value = '0.0₃0774'
$ws.Range('E20').Value = '  -2.84%  '
$ws.Range('E21').Value = '  -2.31%  '
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('E23').Value = '  -1.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.21'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.56%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '170.24'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.82'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.96%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.52'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.54%  '
$ws.Range('E28').Value = '  +1.82%  '
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.23'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.93%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.79'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.50%  '
$ws.Range('E32').Value = '  -2.76%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.85'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.29%  '
$ws.Range('E34').Value = '  -0.31%  '
$ws.Range('D35').Value = '1.362.80'
$ws.Range('E35').Value = '  -2.45%  '
$ws.Range('E36').Value = '  -4.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.36'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.86%  '
$ws.Range('E39').Value = '  -1.82%  '
$ws.Range('B40').Value = 'HuobiToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.45'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.42%  '
$ws.Range('B41').Value = 'WEMIXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.22'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.31%  '
$ws.Range('E42').Value = '  -1.34%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '81.80'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.65%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.80'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.82%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.80'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.15%  '
$ws.Range('E46').Value = '  +1.62%  '
$ws.Range('D47').Value = '1.975.56'
$ws.Range('E47').Value = '  +0.43%  '
$ws.Range('E48').Value = '  -4.57%  '
$ws.Range('E49').Value = '  -0.16%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '102.09'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.83%  '
$ws.Range('E51').Value = '  -4.82%  '
